$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "A2" = 12.32780537112779
    "B2" = 0.0000000000000001110223024625157
    "C2" = 0.005650367507975135
    "D2" = 0.6847753630127387
    "E2" = 0.4689172977892279

    "A3" = 9.706540520275604
    "B3" = 0.0000000000000001110223024625157
    "C3" = 0.005016216216537845
    "D3" = 0.6079217459363185
    "E3" = 0.3695688491822618

    "A4" = 11.53359894141371
    "B4" = 0.0000000000000001110223024625157
    "C4" = 0.005384636403213265
    "D4" = 0.65257106595237
    "E4" = 0.4258489961182125

    "A5" = 11.09587213457653
    "B5" = 0.0000000000000001110223024625157
    "C5" = 0.0058218375788665
    "D5" = 0.7055560431852629
    "E5" = 0.4978093300752446

    "A6" = 9.019377403777067
    "B6" = 0.0000000000000001110223024625157
    "C6" = 0.003979624763272925
    "D6" = 0.4822958839541546
    "E6" = 0.2326093196791194
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
